# Apply updated specific-consumption values (Electricity / Heat) to the
# "glass" sheet, fixing the bug caused by mutable default parameters so
# the rest-sector calculation now produces per-country values instead of
# repeating the same defaults (1.39 / 6.13) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("glass")

$updates = @{
    2  = @(1.58, 7.48)
    3  = @(1.38, 8.34)
    5  = @(2.45, 4.62)
    7  = @(2.97, 8.94)
    8  = @(1.34, 5.71)
    9  = @(1.29, 6.28)
    10 = @(1.44, 6.75)
    11 = @(2.15, 7.83)
    12 = @(1.13, 5.99)
    13 = @(3.38, 10.32)
    15 = @(1.9, 7.02)
    18 = @(1.69, 6.22)
    19 = @(1.1, 5.74)
    20 = @(1.18, 7.16)
    21 = @(3.36, 4.31)
    22 = @(2.36, 7.02)
    23 = @(2.25, 6.4)
    24 = @(1.91, 5.94)
    25 = @(1.28, 6.01)
    31 = @(0.99, 5.41)
    34 = @(3.38, 4.12)
    35 = @(0.99, 5.41)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
